# Auto-generated edit script applying scheduled runner updates to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 167.5
$ws.Range("I33").Value = 167.5
$ws.Range("K33").Value = 167.5
$ws.Range("M33").Value = 61.5
$ws.Range("H69").Value = 33984
$ws.Range("J69").Value = 19979.166
$ws.Range("L69").Value = 59937.49800000001
$ws.Range("N69").Value = -61685.49800000001
$ws.Range("H72").Value = 33984
$ws.Range("J72").Value = 19979.166
$ws.Range("L72").Value = 179812.494
$ws.Range("N72").Value = -188548.494
$ws.Range("H88").Value = 1815.5454
$ws.Range("I88").Value = 1907.8
$ws.Range("J88").Value = 1738.6666
$ws.Range("K88").Value = 1907.8
$ws.Range("L88").Value = 1738.6666
$ws.Range("M88").Value = -1501.8
$ws.Range("N88").Value = -2550.6666
$ws.Range("H91").Value = 1815.5454
$ws.Range("I91").Value = 1907.8
$ws.Range("J91").Value = 1738.6666
$ws.Range("K91").Value = 1907.8
$ws.Range("L91").Value = 1738.6666
$ws.Range("M91").Value = -503.8
$ws.Range("N91").Value = -4546.6666
$ws.Range("H116").Value = 8378.9
$ws.Range("I116").Value = 9332
$ws.Range("K116").Value = 9332
$ws.Range("M116").Value = -5890
$ws.Range("H125").Value = 6347.6924
$ws.Range("J125").Value = 10724.857
$ws.Range("L125").Value = 96523.713
$ws.Range("N125").Value = -101443.713
$ws.Range("H137").Value = 1653.1666
$ws.Range("I137").Value = 1605.625
$ws.Range("K137").Value = 4816.875
$ws.Range("M137").Value = -2266.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 807.6667
$ws.Range("I4").Value = 461.75
$ws.Range("K4").Value = 461.75
$ws.Range("M4").Value = -345.75
$ws.Range("H46").Value = 11498
$ws.Range("J46").Value = 4997
$ws.Range("L46").Value = 4997
$ws.Range("N46").Value = -5635
$ws.Range("H50").Value = 5073.143
$ws.Range("I50").Value = 10116
$ws.Range("K50").Value = 10116
$ws.Range("M50").Value = -9402
$ws.Range("H110").Value = 4626295
$ws.Range("I110").Value = 5286766
$ws.Range("K110").Value = 5286766
$ws.Range("M110").Value = -5284721
$ws.Range("H132").Value = 2133.3333
$ws.Range("I132").Value = 2133.3333
$ws.Range("K132").Value = 6399.999899999999
$ws.Range("M132").Value = -3869.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3188.1667
$ws.Range("I107").Value = 3782.75
$ws.Range("K107").Value = 3782.75
$ws.Range("M107").Value = -1862.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 231749.5
$ws.Range("J3").Value = 9999
$ws.Range("L3").Value = 9999
$ws.Range("N3").Value = -10225
$ws.Range("H5").Value = 568.44446
$ws.Range("I5").Value = 259.83334
$ws.Range("J5").Value = 1185.6666
$ws.Range("K5").Value = 259.83334
$ws.Range("L5").Value = 1185.6666
$ws.Range("M5").Value = -147.83334
$ws.Range("N5").Value = -1409.6666
$ws.Range("H8").Value = 783.6
$ws.Range("I8").Value = 639.6667
$ws.Range("J8").Value = 999.5
$ws.Range("K8").Value = 639.6667
$ws.Range("L8").Value = 999.5
$ws.Range("M8").Value = -499.6667
$ws.Range("N8").Value = -1279.5
$ws.Range("H10").Value = 1850.5
$ws.Range("I10").Value = 2251
$ws.Range("J10").Value = 1450
$ws.Range("K10").Value = 2251
$ws.Range("L10").Value = 1450
$ws.Range("M10").Value = -2112
$ws.Range("N10").Value = -1728
$ws.Range("H13").Value = 2317.3333
$ws.Range("J13").Value = 2317.3333
$ws.Range("L13").Value = 2317.3333
$ws.Range("N13").Value = -2595.3333
$ws.Range("H14").Value = 1875
$ws.Range("J14").Value = 1875
$ws.Range("L14").Value = 1875
$ws.Range("N14").Value = -2215
$ws.Range("H15").Value = 5950
$ws.Range("I15").Value = 9975
$ws.Range("J15").Value = 583.3333
$ws.Range("K15").Value = 9975
$ws.Range("L15").Value = 583.3333
$ws.Range("M15").Value = -9805
$ws.Range("N15").Value = -923.3333
$ws.Range("H17").Value = 1377
$ws.Range("I17").Value = 1169.3334
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 1169.3334
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -995.3334
$ws.Range("N17").Value = -2348
$ws.Range("H19").Value = 106.333336
$ws.Range("I19").Value = 107.125
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 107.125
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = 62.875
$ws.Range("N19").Value = -440
$ws.Range("H22").Value = 321.85715
$ws.Range("I22").Value = 290.8
$ws.Range("J22").Value = 399.5
$ws.Range("K22").Value = 290.8
$ws.Range("L22").Value = 399.5
$ws.Range("M22").Value = 59.19999999999999
$ws.Range("N22").Value = -1099.5
$ws.Range("H24").Value = 106.333336
$ws.Range("I24").Value = 107.125
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 107.125
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = 62.875
$ws.Range("N24").Value = -440
$ws.Range("H25").Value = 602.1667
$ws.Range("I25").Value = 602.1667
$ws.Range("K25").Value = 602.1667
$ws.Range("M25").Value = -428.1667
$ws.Range("H31").Value = 2334.9167
$ws.Range("I31").Value = 2334.9167
$ws.Range("K31").Value = 2334.9167
$ws.Range("M31").Value = -2039.9167
$ws.Range("H34").Value = 2334.9167
$ws.Range("I34").Value = 2334.9167
$ws.Range("K34").Value = 2334.9167
$ws.Range("M34").Value = -2132.9167
$ws.Range("H74").Value = 39166.668
$ws.Range("J74").Value = 39166.668
$ws.Range("L74").Value = 39166.668
$ws.Range("N74").Value = -40914.668
$ws.Range("H77").Value = 39166.668
$ws.Range("J77").Value = 39166.668
$ws.Range("L77").Value = 117500.004
$ws.Range("N77").Value = -126236.004
$ws.Range("H86").Value = 5707.273
$ws.Range("I86").Value = 5694
$ws.Range("K86").Value = 5694
$ws.Range("M86").Value = -4571
$ws.Range("H88").Value = 17340
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 17340
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 17340
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -18152
$ws.Range("H89").Value = 5707.273
$ws.Range("I89").Value = 5694
$ws.Range("K89").Value = 28470
$ws.Range("M89").Value = -22854
$ws.Range("H91").Value = 17340
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 17340
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 17340
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -20148
$ws.Range("H99").Value = 1899.9166
$ws.Range("J99").Value = 1899.5
$ws.Range("L99").Value = 1899.5
$ws.Range("N99").Value = -4895.5
$ws.Range("H126").Value = 1899.9166
$ws.Range("J126").Value = 1899.5
$ws.Range("L126").Value = 5698.5
$ws.Range("N126").Value = -10638.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3417.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3417.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10252.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -10606.5
$ws.Range("H86").Value = 732.6667
$ws.Range("J86").Value = 732.6667
$ws.Range("L86").Value = 2198.0001
$ws.Range("N86").Value = -4570.0001
$ws.Range("H89").Value = 732.6667
$ws.Range("J89").Value = 732.6667
$ws.Range("L89").Value = 6594.0003
$ws.Range("N89").Value = -18450.0003
$ws.Range("H115").Value = 3975
$ws.Range("J115").Value = 3975
$ws.Range("L115").Value = 11925
$ws.Range("N115").Value = -14275

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H69").Value = 49999
$ws.Range("J69").Value = 49999
$ws.Range("L69").Value = 49999
$ws.Range("N69").Value = -51497
$ws.Range("H72").Value = 49999
$ws.Range("J72").Value = 49999
$ws.Range("L72").Value = 149997
$ws.Range("N72").Value = -157485
$ws.Range("H94").Value = 25200.928
$ws.Range("J94").Value = 27235.637
$ws.Range("L94").Value = 27235.637
$ws.Range("N94").Value = -28587.637
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 6075.7144
$ws.Range("I122").Value = 5422
$ws.Range("K122").Value = 16266
$ws.Range("M122").Value = -13816
$ws.Range("H132").Value = 16900
$ws.Range("I132").Value = 16900
$ws.Range("K132").Value = 50700
$ws.Range("M132").Value = -48170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1463.091
$ws.Range("I136").Value = 1075.125
$ws.Range("K136").Value = 3225.375
$ws.Range("M136").Value = -675.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 19759.6
$ws.Range("J45").Value = 19759.6
$ws.Range("L45").Value = 19759.6
$ws.Range("N45").Value = -20741.6
$ws.Range("H122").Value = 1492.8572
$ws.Range("I122").Value = 1492.8572
$ws.Range("K122").Value = 4478.571599999999
$ws.Range("M122").Value = -2028.571599999999
